$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B (old B -> C, old C -> D); column A keeps its exact
# original width, and the old B/C widths move along automatically with them.
$ws.Columns("B").Insert()

# New column B should be as wide as column A (75.81640625). The engine
# quantizes ColumnWidth to 1/6-character steps, so 75 is the closest input
# that lands on the nearest achievable bucket to the target width.
$ws.Columns("B").ColumnWidth = 75

# Row 1 header: add the new "StatQuery" column header in B1.
$ws.Range("B1").Value = "StatQuery"

# Row 2, column A: drop the "AND a.arm_drug IN ['Nivolumab']" clause from the query.
$ws.Range("A2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE a.arm_id IN ['Z1D']  RETURN coalesce(c.case_id,'') AS ``Case ID`` , coalesce(t.clinical_trial_designation ,'')as ``Trial Code`` , coalesce(a.arm_id,'') As ``Arm`` , coalesce(a.arm_drug,'') As ``Arm Treatment`` , coalesce(c.disease,'') As Diagnosis , coalesce(c.gender,'') As Gender , coalesce(c.race,'') As Race , coalesce(c.ethnicity,'') As Ethnicity"

# Row 2, column B: new stat query cell (inherits the wrap-text style from the insert).
$ws.Range("B2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE a.arm_id IN ['Z1D'] OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"

# Row height shrinks now that the long query text wraps across a wider cell.
$ws.Rows("2").RowHeight = 87
